$d = $word.ActiveDocument

# --- Section: "Aim to get associated..." career aspiration paragraph ---

# Step 1: " associated with a progressive organization that provides me opportunities"
#         -> " associated with a progressive organization that provides me opportunit"
#         (stays entirely inside the original run; keeps its rPr: rFonts ascii=Cambria hAnsi=Cambria)
$d.Content.Find.Execute(
    " associated with a progressive organization that provides me opportunities",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " associated with a progressive organization that provides me opportunit", 2) | Out-Null

# Step 2: the single space run between "...opportunit" and "to work in JavaScript..."
#         becomes "ies " (keeps the original run's rPr: rFonts hAnsi=Cambria only)
$anchor = $d.Content
$anchor.Find.Execute("to work in JavaScript projects where I can", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$gap = $anchor.Duplicate
$gap.Collapse(1) | Out-Null
$gap.MoveStart(1, -1) | Out-Null
$gap.Text = "ies "

# Step 3: " sharpen my skills" -> " Sharpen" (same run/rPr) + new " my skills" run (same rPr: hAnsi=Cambria only)
$r = $d.Content
$r.Find.Execute(" sharpen my skills", $true, $false, $false, $false, $false,
                $true, 1, $false, "", 0) | Out-Null
$r.Text = " Sharpen"
$tail = $r.Duplicate
$tail.Collapse(0) | Out-Null
$tail.InsertAfter(" my skills") | Out-Null

# Step 4: "in customer satisfaction " ->
#   "in customer and my managers’ satisfaction and focus on producing a quality code "
#   (all new pieces share the original run's rPr: rFonts ascii=Cambria hAnsi=Cambria, so a
#    single in-place text replacement reproduces the same final formatting)
$d.Content.Find.Execute(
    "in customer satisfaction ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "in customer and my managers’ satisfaction and focus on producing a quality code ", 2) | Out-Null

# Step 5: "my managers and peers. I am " -> "and peers. I am " (same run/rPr)
$d.Content.Find.Execute(
    "my managers and peers. I am ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and peers. I am ", 2) | Out-Null
